$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the "binary of certain length), population" text
# (previously split across runs by a grammar-check proofErr pair) into
# a single clean sentence, also fixing the double space after the
# comma ("),  population" -> "), population").
# ---------------------------------------------------------------------
$old1 = "first specify how we describe sequences (binary of certain length" + `
        "),  population of sequences, analogous to figure 1"
$new1 = "first specify how we describe sequences (binary of certain length" + `
        "), population of sequences, analogous to figure 1"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Could not find text for change 1"
}

# ---------------------------------------------------------------------
# Change 2: turn "discussion: M didn't think about it too carefully"
# into "discussion: (M didn't think about it too carefully yet)", and
# relocate the (hidden) "_GoBack" bookmark so it now sits right before
# the new "yet)" text instead of after "... has a bit too much
# attention" in the following bullet.
# ---------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("discussion: M", $true, $false, $false, $false, `
                                   $false, $true, 1, $false, "discussion: (M", 2)
if (-not $found2) {
    throw "Could not find 'discussion: M'"
}

$sentinel = $d.Content
$found3 = $sentinel.Find.Execute("too carefully")
if (-not $found3) {
    throw "Could not find 'too carefully'"
}
$insertPoint = $sentinel.End

# Insert the trailing " yet)" right after "...too carefully".
$spaceRange = $d.Range($insertPoint, $insertPoint)
$spaceRange.InsertAfter(" ")

$bookmarkPos = $insertPoint + 1
$yetRange = $d.Range($bookmarkPos, $bookmarkPos)
$yetRange.InsertAfter("yet)")

# Re-adding a bookmark with the same name moves it from its previous
# location (exactly like real Word does), so this relocates "_GoBack"
# from the "MPL paper by Barton..." bullet to right before "yet)".
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
